$d = $word.ActiveDocument

# 1. Center the title paragraph (first paragraph in the document).
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Alignment = 1  # wdAlignParagraphCenter

# 2. Remove the old "_GoBack" bookmark (sits after "Usuário consegue criar, excluir e marcar tarefas.")
#    and re-add it collapsed at the very start of the title paragraph.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
# A bookmark range collapsed exactly at document position 0 gets mis-handled
# (it balloons to cover the whole first paragraph), so nudge a throwaway
# character in front of the title first, anchor the bookmark just after it,
# then delete the throwaway character; the bookmark collapses back to 0-0.
$d.Range(0, 0).InsertBefore("X")
$bmRange = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $bmRange)
$d.Range(0, 1).Delete()

# 3. Delete the whole paragraph "Somente o dono da conta pode acessar suas tarefas."
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Somente o dono da conta pode acessar suas tarefas.*") {
        $p.Range.Delete()
        break
    }
}

# 4. Drop the stale <w:lastRenderedPageBreak/> cached before "8. Responsabilidades" by
#    re-typing the run's text through Find/Replace (forces Word to rebuild the run).
$d.Content.Find.Execute("8. Responsabilidades", $true, $false, $false, $false, $false,
                         $true, 1, $false, "8. Responsabilidades", 2) | Out-Null

# 5. Shrink the section page margins to 720 twips (0.5") on every side.
#    PageSetup margins are expressed in points (1 pt = 20 twips), so 720
#    twips == 36 points.
$section = $d.Sections.Item(1)
$section.PageSetup.TopMargin = 36
$section.PageSetup.BottomMargin = 36
$section.PageSetup.LeftMargin = 36
$section.PageSetup.RightMargin = 36
